# Auto-generated: updates Coin/Link/Price/Volume(1h) cells per the symbol-list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''305.78'
$ws.Range("E2").Value = '''1.03%'
$ws.Range("D3").Value = '''36.03'
$ws.Range("E3").Value = '''-1.44%'
$ws.Range("D4").Value = '''5.067'
$ws.Range("E4").Value = '''1.64%'
$ws.Range("D5").Value = '''0.08004'
$ws.Range("E5").Value = '''3.13%'
$ws.Range("D6").Value = '''2.177'
$ws.Range("E6").Value = '''3.60%'
$ws.Range("D7").Value = '''8.017'
$ws.Range("E7").Value = '''1.33%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9291'
$ws.Range("E8").Value = '''0.79%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.09894'
$ws.Range("E9").Value = '''0.97%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1867'
$ws.Range("E10").Value = '''0.03%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09029'
$ws.Range("E11").Value = '''5.18%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03628'
$ws.Range("E12").Value = '''3.44%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09926'
$ws.Range("E13").Value = '''-0.22%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001453'
$ws.Range("E14").Value = '''-0.72%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005665'
$ws.Range("E15").Value = '''0.89%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.453'
$ws.Range("E16").Value = '''-0.33%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.153'
$ws.Range("E17").Value = '''2.74%'
$ws.Range("E18").Value = '''13.87%'
$ws.Range("E19").Value = '''-1.08%'
$ws.Range("D20").Value = '''0.1357'
$ws.Range("E20").Value = '''1.01%'
$ws.Range("D21").Value = '''5.066'
$ws.Range("E21").Value = '''6.57%'
$ws.Range("E22").Value = '''-0.23%'
$ws.Range("D23").Value = '''0.04590'
$ws.Range("E23").Value = '''0.31%'
$ws.Range("E24").Value = '''0.92%'
$ws.Range("D25").Value = '''0.004750'
$ws.Range("E25").Value = '''-6.51%'
$ws.Range("D26").Value = '''0.0001301'
$ws.Range("E26").Value = '''-6.85%'
$ws.Range("D27").Value = '''0.0004507'
$ws.Range("E27").Value = '''65.34%'
$ws.Range("D39").Value = '''0.01946'
$ws.Range("E39").Value = '''10.34%'
$ws.Range("D40").Value = '''0.04895'
$ws.Range("E40").Value = '''4.88%'
$ws.Range("D41").Value = '''0.007813'
$ws.Range("E41").Value = '''4.80%'
$ws.Range("D42").Value = '''0.1395'
$ws.Range("E42").Value = '''0.50%'
$ws.Range("D43").Value = '''0.007818'
$ws.Range("E43").Value = '''1.56%'
$ws.Range("D44").Value = '''0.002106'
$ws.Range("E44").Value = '''-5.76%'
$ws.Range("D45").Value = '''0.01142'
$ws.Range("E45").Value = '''9.60%'
$ws.Range("D46").Value = '''0.00006316'
$ws.Range("E46").Value = '''2.28%'
$ws.Range("D47").Value = '''0.00000000751'
$ws.Range("E47").Value = '''0.32%'
$ws.Range("D48").Value = '''52.03'
$ws.Range("E48").Value = '''36.34%'
$ws.Range("D49").Value = '''0.001803'
$ws.Range("E49").Value = '''-9.65%'
$ws.Range("D50").Value = '''0.00002102'
$ws.Range("E50").Value = '''0.32%'
$ws.Range("D51").Value = '''0.0002002'
$ws.Range("E51").Value = '''0.32%'
